$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Min values that changed
$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 4.2

# Remove the "theta_threshold_range" row (row 5); rows below shift up
$ws.Rows.Item(5).Delete()

# Update the value that moved up from the old row 6 ("pie_threshold_range")
$ws.Range("C5").Value = 15

# Normalize the two cells that carried the special (Times New Roman) font
# so they match the plain data style, the same way the rest of the data
# cells are styled.
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Move the active selection to C4, matching the saved workbook state
$ws.Range("C4").Select() | Out-Null
